$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Thbs1"
$ws.Cells.Item(2, 3).Value = "Itga3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 142.9073533333333
$ws.Cells.Item(2, 8).Value = 428.72206
$ws.Cells.Item(2, 9).Value = 0.5576664151504187
$ws.Cells.Item(2, 10).Value = 0.5576664151504188
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.269168666666666
$ws.Cells.Item(2, 14).Value = 24.807506
$ws.Cells.Item(2, 15).Value = 0.671680253471746
$ws.Cells.Item(2, 16).Value = 0.671680253471746
$ws.Cells.Item(2, 17).Value = 1181.725008420262
$ws.Cells.Item(2, 18).Value = 10635.52507578236
$ws.Cells.Item(2, 19).Value = 0.3745735190809131
$ws.Cells.Item(2, 20).Value = 0.3745735190809132

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Thbs1"
$ws.Cells.Item(3, 3).Value = "Itga3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 142.9073533333333
$ws.Cells.Item(3, 8).Value = 428.72206
$ws.Cells.Item(3, 9).Value = 0.5576664151504187
$ws.Cells.Item(3, 10).Value = 0.5576664151504188
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.121294
$ws.Cells.Item(3, 14).Value = 0.363882
$ws.Cells.Item(3, 15).Value = 0.009852354928133683
$ws.Cells.Item(3, 16).Value = 0.009852354928133683
$ws.Cells.Item(3, 17).Value = 17.33380451521333
$ws.Cells.Item(3, 18).Value = 156.00424063692
$ws.Cells.Item(3, 19).Value = 0.005494327453561872
$ws.Cells.Item(3, 20).Value = 0.005494327453561873

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Thbs1"
$ws.Cells.Item(4, 3).Value = "Itga3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 142.9073533333333
$ws.Cells.Item(4, 8).Value = 428.72206
$ws.Cells.Item(4, 9).Value = 0.5576664151504187
$ws.Cells.Item(4, 10).Value = 0.5576664151504188
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.920705666666667
$ws.Cells.Item(4, 14).Value = 11.762117
$ws.Cells.Item(4, 15).Value = 0.3184673916001203
$ws.Cells.Item(4, 16).Value = 0.3184673916001203
$ws.Cells.Item(4, 17).Value = 560.2976700223355
$ws.Cells.Item(4, 18).Value = 5042.67903020102
$ws.Cells.Item(4, 19).Value = 0.1775985686159436
$ws.Cells.Item(4, 20).Value = 0.1775985686159437

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Thbs1"
$ws.Cells.Item(5, 3).Value = "Itga3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 63.967809
$ws.Cells.Item(5, 8).Value = 191.903427
$ws.Cells.Item(5, 9).Value = 0.2496211559306514
$ws.Cells.Item(5, 10).Value = 0.2496211559306514
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 8.269168666666666
$ws.Cells.Item(5, 14).Value = 24.807506
$ws.Cells.Item(5, 15).Value = 0.671680253471746
$ws.Cells.Item(5, 16).Value = 0.671680253471746
$ws.Cells.Item(5, 17).Value = 528.9606018581179
$ws.Cells.Item(5, 18).Value = 4760.645416723061
$ws.Cells.Item(5, 19).Value = 0.1676656012874101
$ws.Cells.Item(5, 20).Value = 0.1676656012874101

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Thbs1"
$ws.Cells.Item(6, 3).Value = "Itga3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 63.967809
$ws.Cells.Item(6, 8).Value = 191.903427
$ws.Cells.Item(6, 9).Value = 0.2496211559306514
$ws.Cells.Item(6, 10).Value = 0.2496211559306514
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.121294
$ws.Cells.Item(6, 14).Value = 0.363882
$ws.Cells.Item(6, 15).Value = 0.009852354928133683
$ws.Cells.Item(6, 16).Value = 0.009852354928133683
$ws.Cells.Item(6, 17).Value = 7.758911424846
$ws.Cells.Item(6, 18).Value = 69.83020282361399
$ws.Cells.Item(6, 19).Value = 0.00245935622579978
$ws.Cells.Item(6, 20).Value = 0.00245935622579978

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Thbs1"
$ws.Cells.Item(7, 3).Value = "Itga3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 63.967809
$ws.Cells.Item(7, 8).Value = 191.903427
$ws.Cells.Item(7, 9).Value = 0.2496211559306514
$ws.Cells.Item(7, 10).Value = 0.2496211559306514
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.920705666666667
$ws.Cells.Item(7, 14).Value = 11.762117
$ws.Cells.Item(7, 15).Value = 0.3184673916001203
$ws.Cells.Item(7, 16).Value = 0.3184673916001203
$ws.Cells.Item(7, 17).Value = 250.798951230551
$ws.Cells.Item(7, 18).Value = 2257.190561074959
$ws.Cells.Item(7, 19).Value = 0.07949619841744143
$ws.Cells.Item(7, 20).Value = 0.07949619841744145

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Thbs1"
$ws.Cells.Item(8, 3).Value = "Itga3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 49.38440333333333
$ws.Cells.Item(8, 8).Value = 148.15321
$ws.Cells.Item(8, 9).Value = 0.1927124289189298
$ws.Cells.Item(8, 10).Value = 0.1927124289189298
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.269168666666666
$ws.Cells.Item(8, 14).Value = 24.807506
$ws.Cells.Item(8, 15).Value = 0.671680253471746
$ws.Cells.Item(8, 16).Value = 0.671680253471746
$ws.Cells.Item(8, 17).Value = 408.3679606660288
$ws.Cells.Item(8, 18).Value = 3675.311645994259
$ws.Cells.Item(8, 19).Value = 0.1294411331034226
$ws.Cells.Item(8, 20).Value = 0.1294411331034226

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Thbs1"
$ws.Cells.Item(9, 3).Value = "Itga3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 49.38440333333333
$ws.Cells.Item(9, 8).Value = 148.15321
$ws.Cells.Item(9, 9).Value = 0.1927124289189298
$ws.Cells.Item(9, 10).Value = 0.1927124289189298
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.121294
$ws.Cells.Item(9, 14).Value = 0.363882
$ws.Cells.Item(9, 15).Value = 0.009852354928133683
$ws.Cells.Item(9, 16).Value = 0.009852354928133683
$ws.Cells.Item(9, 17).Value = 5.990031817913333
$ws.Cells.Item(9, 18).Value = 53.91028636122
$ws.Cells.Item(9, 19).Value = 0.00189867124877203
$ws.Cells.Item(9, 20).Value = 0.001898671248772031

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Thbs1"
$ws.Cells.Item(10, 3).Value = "Itga3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 49.38440333333333
$ws.Cells.Item(10, 8).Value = 148.15321
$ws.Cells.Item(10, 9).Value = 0.1927124289189298
$ws.Cells.Item(10, 10).Value = 0.1927124289189298
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.920705666666667
$ws.Cells.Item(10, 14).Value = 11.762117
$ws.Cells.Item(10, 15).Value = 0.3184673916001203
$ws.Cells.Item(10, 16).Value = 0.3184673916001203
$ws.Cells.Item(10, 17).Value = 193.6217099939522
$ws.Cells.Item(10, 18).Value = 1742.59538994557
$ws.Cells.Item(10, 19).Value = 0.06137262456673517
$ws.Cells.Item(10, 20).Value = 0.06137262456673517